$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates per the refreshed cryptos snapshot.
# Values that look numeric are prefixed with a literal apostrophe so Excel
# keeps them as text (matching the original text-typed cells) instead of
# coercing them into numbers (which would drop formatting like trailing zeros).

$ws.Range("D2").Value = "62.974.52"
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").Value = "3.224.60"
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'527.41"
$ws.Range("E5").Value = "  +3.77%  "
$ws.Range("D6").Value = "'170.15"
$ws.Range("E6").Value = "  -2.44%  "
$ws.Range("D7").Value = "'0.595"
$ws.Range("E7").Value = "  +2.03%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "3.221.25"
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("E10").Value = "  -0.26%  "
$ws.Range("D11").Value = "'52.79"
$ws.Range("E11").Value = "  -6.48%  "
$ws.Range("D12").Value = "'0.134"
$ws.Range("E12").Value = "  +4.29%  "
$ws.Range("D13").Value = "'0.0000253"
$ws.Range("E13").Value = "  +1.31%  "
$ws.Range("D14").Value = "'9.11"
$ws.Range("E14").Value = "  +2.16%  "
$ws.Range("D15").Value = "3.738.94"
$ws.Range("E15").Value = "  -0.50%  "
$ws.Range("E16").Value = "  -0.79%  "
$ws.Range("D17").Value = "3.228.81"
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("D18").Value = "62.846.41"
$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("D19").Value = "'17.15"
$ws.Range("E19").Value = "  +1.46%  "
$ws.Range("D20").Value = "'11.04"
$ws.Range("E20").Value = "  +4.26%  "
$ws.Range("E21").Value = "  +4.17%  "
$ws.Range("D22").Value = "'365.45"
$ws.Range("E22").Value = "  +0.53%  "
$ws.Range("D23").Value = "'3.74"
$ws.Range("E23").Value = "  +4.78%  "
$ws.Range("D24").Value = "'81.02"
$ws.Range("E24").Value = "  +3.18%  "
$ws.Range("D25").Value = "'11.13"
$ws.Range("E25").Value = "  +3.81%  "
$ws.Range("D26").Value = "'3.95"
$ws.Range("E26").Value = "  +6.49%  "
$ws.Range("D27").Value = "'6.06"
$ws.Range("E27").Value = "  +0.86%  "
$ws.Range("E28").Value = "  +1.51%  "
$ws.Range("D29").Value = "'11.23"
$ws.Range("E29").Value = "  +1.88%  "
$ws.Range("D30").Value = "'8.18"
$ws.Range("E30").Value = "  +0.39%  "
$ws.Range("D31").Value = "'28.38"
$ws.Range("E31").Value = "  +1.81%  "
$ws.Range("D32").Value = "'631.60"
$ws.Range("E32").Value = "  -0.71%  "
$ws.Range("E33").Value = "  -1.71%  "
$ws.Range("D34").Value = "'11.16"
$ws.Range("E34").Value = "  +2.74%  "
$ws.Range("D35").Value = "'0.105"
$ws.Range("E35").Value = "  +4.31%  "
$ws.Range("D36").Value = "'56.64"
$ws.Range("E36").Value = "  -3.15%  "
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("D38").Value = "'36.55"
$ws.Range("E38").Value = "  +4.01%  "
$ws.Range("E39").Value = "  +1.70%  "
$ws.Range("E40").Value = "  +0.16%  "
$ws.Range("D41").Value = "0.0₃0709"
$ws.Range("E41").Value = "  +12.47%  "
$ws.Range("D42").Value = "'0.122"
$ws.Range("E42").Value = "  +2.21%  "
$ws.Range("E43").Value = "  +10.91%  "
$ws.Range("D44").Value = "2.866.34"
$ws.Range("E44").Value = "  +1.84%  "
$ws.Range("D45").Value = "'2.93"
$ws.Range("E45").Value = "  +6.99%  "
$ws.Range("D46").Value = "'2.67"
$ws.Range("E46").Value = "  +3.69%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "'3.07"
$ws.Range("E47").Value = "  +6.96%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "'0.0391"
$ws.Range("E48").Value = "  +4.62%  "
$ws.Range("E49").Value = "  -1.27%  "
$ws.Range("E50").Value = "  +2.77%  "
$ws.Range("D51").Value = "'133.83"
$ws.Range("E51").Value = "  +2.20%  "
